# Weekly update: insert a new price record as row 23 (for the week of
# 2022-08-26) on the "Fruta/Guayaba" sheet. All existing data rows from
# 23 downward shift down by one row to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23:50 down to 24:51, inserting a blank row 23.
$ws.Rows.Item(23).Insert()

# Seed the new row 23 with the same constant/static columns as the row
# immediately below it (now row 24, originally row 23), since every
# record on this sheet shares the same market/product/unit metadata.
$ws.Range("A23:T23").Value = $ws.Range("A24:T24").Value()

# Now overwrite the fields that are specific to this new weekly record.
$ws.Range("D23").Value = 44799      # Fecha: 2022-08-26
$ws.Range("L23").Value = "Primera"  # Calidad
$ws.Range("M23").Value = 200        # Volumen
$ws.Range("N23").Value = 600        # Precio minimo
$ws.Range("O23").Value = 700        # Precio maximo
$ws.Range("P23").Value = 650        # Precio promedio ponderado
$ws.Range("S23").Value = 650        # Precio $/Kg
